{"js": "// Replace each \"NNN\u00d7N=\" expression in the practice-table cells with its\n// updated value (problem set refreshed to a new day's numbers). Every old\n// value is unique in the document, so an exact, case-sensitive search for\n// each old string -> insertText(..., \"Replace\") on the single hit is\n// unambiguous and leaves everything else (the date line, empty answer\n// rows, formatting) untouched.\nconst replacements = [\n  [\"514\u00d74=\", \"422\u00d78=\"],\n  [\"566\u00d76=\", \"848\u00d76=\"],\n  [\"164\u00d76=\", \"536\u00d78=\"],\n  [\"778\u00d79=\", \"777\u00d73=\"],\n  [\"830\u00d72=\", \"989\u00d78=\"],\n  [\"275\u00d79=\", \"978\u00d76=\"],\n  [\"413\u00d72=\", \"161\u00d78=\"],\n  [\"479\u00d77=\", \"639\u00d74=\"],\n  [\"721\u00d79=\", \"105\u00d78=\"],\n  [\"275\u00d76=\", \"280\u00d78=\"],\n  [\"533\u00d77=\", \"953\u00d73=\"],\n  [\"785\u00d74=\", \"699\u00d78=\"],\n  [\"333\u00d72=\", \"523\u00d79=\"],\n  [\"348\u00d74=\", \"878\u00d79=\"],\n  [\"319\u00d79=\", \"200\u00d79=\"],\n  [\"277\u00d77=\", \"506\u00d75=\"],\n  [\"430\u00d78=\", \"966\u00d74=\"],\n  [\"480\u00d75=\", \"584\u00d75=\"],\n  [\"217\u00d76=\", \"753\u00d79=\"],\n  [\"664\u00d79=\", \"993\u00d74=\"],\n  [\"389\u00d77=\", \"964\u00d73=\"],\n  [\"547\u00d78=\", \"534\u00d75=\"],\n  [\"378\u00d79=\", \"413\u00d74=\"],\n  [\"616\u00d73=\", \"396\u00d74=\"],\n  [\"699\u00d73=\", \"837\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Expected to find \"${oldText}\" exactly once, but found none.`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Refresh the practice-table \"NNN\u00d7N=\" expressions to the new day's numbers.\n# Every old value is unique in the document, so a straight Find/Replace\n# (one pair at a time, ReplaceAll just to be safe) is unambiguous and\n# leaves everything else (the date line, blank answer rows, formatting)\n# untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"514\u00d74=\", \"422\u00d78=\"),\n    @(\"566\u00d76=\", \"848\u00d76=\"),\n    @(\"164\u00d76=\", \"536\u00d78=\"),\n    @(\"778\u00d79=\", \"777\u00d73=\"),\n    @(\"830\u00d72=\", \"989\u00d78=\"),\n    @(\"275\u00d79=\", \"978\u00d76=\"),\n    @(\"413\u00d72=\", \"161\u00d78=\"),\n    @(\"479\u00d77=\", \"639\u00d74=\"),\n    @(\"721\u00d79=\", \"105\u00d78=\"),\n    @(\"275\u00d76=\", \"280\u00d78=\"),\n    @(\"533\u00d77=\", \"953\u00d73=\"),\n    @(\"785\u00d74=\", \"699\u00d78=\"),\n    @(\"333\u00d72=\", \"523\u00d79=\"),\n    @(\"348\u00d74=\", \"878\u00d79=\"),\n    @(\"319\u00d79=\", \"200\u00d79=\"),\n    @(\"277\u00d77=\", \"506\u00d75=\"),\n    @(\"430\u00d78=\", \"966\u00d74=\"),\n    @(\"480\u00d75=\", \"584\u00d75=\"),\n    @(\"217\u00d76=\", \"753\u00d79=\"),\n    @(\"664\u00d79=\", \"993\u00d74=\"),\n    @(\"389\u00d77=\", \"964\u00d73=\"),\n    @(\"547\u00d78=\", \"534\u00d75=\"),\n    @(\"378\u00d79=\", \"413\u00d74=\"),\n    @(\"616\u00d73=\", \"396\u00d74=\"),\n    @(\"699\u00d73=\", \"837\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
